$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the existing hyperlink on B2 before clearing the old test data,
# then wipe the two data rows (A2:C2) clean (content + formatting).
$ws.Range("B2").Hyperlinks.Delete()
$ws.Range("A2:C10").Clear()

# Header row (A1:C1) is unchanged - leave as-is.

# New test data rows
$data = @(
    @("John",  "john@test.com.au",  "Nice Service"),
    @("Smith", "smith@test.com.au", "Superb"),
    @("Tony",  "tony@test.com.au",  "So much appreciate your service"),
    @("Brian", "brian@test.com.au", "Excelent Support"),
    @("James", "james@test.com.au", "Keep it up !")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 2), "mailto:" + $row[1])
    $ws.Cells.Item($r, 2).Style = "Hyperlink"
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Widen column C to fit the longer message text.
$ws.Columns("C").ColumnWidth = 29.85

# Update the selected cell shown when the sheet is opened.
$ws.Range("D17").Select()
